$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

$data = @(
    @(2,  "04/10-12:54", 100),
    @(3,  "04/10-12:55", 100),
    @(4,  "04/10-12:57", 100),
    @(5,  "04/10-12:58", 100),
    @(6,  "04/10-13:00", 100),
    @(7,  "04/10-13:01", 100),
    @(8,  "04/10-13:03", 100),
    @(9,  "04/10-13:04", 100),
    @(10, "04/10-13:06", 100),
    @(11, "04/10-13:07", 100),
    @(12, "04/10-13:08", 100),
    @(13, "04/10-13:10", 100),
    @(14, "04/10-13:11", 100),
    @(15, "04/10-13:13", 100),
    @(16, "04/10-13:14", 100),
    @(17, "04/10-13:16", 100),
    @(18, "04/10-13:17", 100),
    @(19, "04/10-13:19", 100),
    @(20, "04/10-13:20", 100),
    @(21, "04/10-13:22", 100),
    @(22, "04/10-13:23", 100),
    @(23, "04/10-13:24", 100),
    @(24, "04/10-13:26", 100),
    @(25, "04/10-13:27", 100),
    @(26, "04/10-13:29", 100),
    @(27, "04/10-13:30", 100),
    @(28, "04/10-13:32", 100),
    @(29, "04/10-13:33", 100),
    @(30, "04/10-13:35", 100)
)

foreach ($row in $data) {
    $r = $row[0]
    $dt = $row[1]
    $val = $row[2]
    $ws.Cells.Item($r, 1).Value = $dt
    $ws.Cells.Item($r, 2).Value = $val
}
